$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log row (row 8) to the sheet, mirroring the existing rows.
$ws.Range("A8").Value = "edit1"
$ws.Range("B8").Value = "riya-morankar"
$ws.Range("C8").Value = "Merged"
# Leading apostrophe forces this date-looking value to be stored as text,
# matching how the Date column is already stored for the other rows.
$ws.Range("E8").Value = "'2025-06-18"
$ws.Range("F8").Value = "N/A"
